$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell to Text format first so the literal string is
# preserved exactly (matches original inline-string cells: prices like
# "261.39" and percentages like "1.70%" must stay as text, not be
# auto-converted to a number/percentage by Excel).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.70%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.759"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.45%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.94%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.641"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8613"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.95%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9245"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.29%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1406"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.53%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05034"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.03%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.18%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03043"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.65%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09089"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.25%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001538"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.54%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006062"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.54%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006127"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.79%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.12%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.166"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.68%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.29%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.54%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.24%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.108"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04262"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.12%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.19%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.78%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.01%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03882"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.60%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.06%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004124"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-34.65%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01502"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.84%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.63%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005305"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.35%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-47.51%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
